$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text, $donor) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

function Set-NumberCell($addr, $num, $donor) {
    $ws.Range($addr).Value = $num
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# --- Cell type / style changes (text <-> number) ---
Set-NumberCell "C14" 1 "C23"
Set-NumberCell "F14" 1 "F29"
Set-NumberCell "C22" 1 "C23"
Set-NumberCell "D22" 1 "D28"
Set-NumberCell "E22" 0 "E15"
Set-TextCell "C28" "0" "C31"
Set-TextCell "C29" "0" "C31"
Set-NumberCell "D29" 1 "D28"
Set-NumberCell "E29" -100 "E15"
Set-TextCell "C30" "0" "C31"
Set-NumberCell "D30" 1 "D28"
Set-NumberCell "E30" -100 "E15"
Set-TextCell "G31" "0" "G14"
Set-TextCell "H31" "***.*" "H14"

# --- Value-only changes ---
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = -63.636363636363
$ws.Range("N14").Value = -91.111111111111
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = -10
$ws.Range("I15").Value = 25
$ws.Range("J15").Value = 27
$ws.Range("K15").Value = -7.407407407407
$ws.Range("L15").Value = 13.636363636363
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -37.5
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = -13.333333333333
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 47
$ws.Range("H16").Value = -29.787234042553
$ws.Range("I16").Value = 149
$ws.Range("J16").Value = 238
$ws.Range("K16").Value = -37.394957983193
$ws.Range("L16").Value = -36.595744680851
$ws.Range("M16").Value = -41.568627450980
$ws.Range("N16").Value = -85.320197044335
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -22.222222222222
$ws.Range("G17").Value = 82
$ws.Range("H17").Value = -3.658536585365
$ws.Range("I17").Value = 349
$ws.Range("J17").Value = 382
$ws.Range("K17").Value = -8.638743455497
$ws.Range("L17").Value = -3.055555555555
$ws.Range("M17").Value = 34.230769230769
$ws.Range("N17").Value = -30.478087649402
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -22.727272727272
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = -19.587628865979
$ws.Range("L18").Value = -46.575342465753
$ws.Range("N18").Value = -87.638668779714
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 190.909090909091
$ws.Range("F19").Value = 72
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = 24.137931034482
$ws.Range("I19").Value = 313
$ws.Range("J19").Value = 328
$ws.Range("K19").Value = -4.573170731707
$ws.Range("L19").Value = -18.701298701298
$ws.Range("M19").Value = 36.681222707423
$ws.Range("N19").Value = 5.387205387205
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 20
$ws.Range("E20").Value = -90
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 56
$ws.Range("H20").Value = -39.285714285714
$ws.Range("I20").Value = 141
$ws.Range("J20").Value = 209
$ws.Range("K20").Value = -32.535885167464
$ws.Range("L20").Value = -28.426395939086
$ws.Range("M20").Value = 62.068965517241
$ws.Range("N20").Value = -86.257309941520
$ws.Range("C21").Value = 69
$ws.Range("D21").Value = 69
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 245
$ws.Range("G21").Value = 275
$ws.Range("H21").Value = -10.909090909090
$ws.Range("I21").Value = 1059
$ws.Range("J21").Value = 1285
$ws.Range("K21").Value = -17.587548638132
$ws.Range("L21").Value = -21.729490022172
$ws.Range("M21").Value = 3.519061583577
$ws.Range("N21").Value = -70.219347581552
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = 44.444444444444
$ws.Range("L22").Value = -35
$ws.Range("M22").Value = -45.833333333333
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 150
$ws.Range("F23").Value = 32
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = 39.130434782608
$ws.Range("I23").Value = 110
$ws.Range("J23").Value = 101
$ws.Range("K23").Value = 8.910891089108
$ws.Range("L23").Value = -16.030534351145
$ws.Range("M23").Value = 80.327868852459
$ws.Range("C24").Value = 51
$ws.Range("D24").Value = 55
$ws.Range("E24").Value = -7.272727272727
$ws.Range("F24").Value = 236
$ws.Range("G24").Value = 197
$ws.Range("H24").Value = 19.796954314720
$ws.Range("I24").Value = 957
$ws.Range("J24").Value = 818
$ws.Range("K24").Value = 16.992665036674
$ws.Range("L24").Value = 15.719467956469
$ws.Range("M24").Value = 84.749034749034
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -21.739130434782
$ws.Range("F25").Value = 94
$ws.Range("G25").Value = 75
$ws.Range("H25").Value = 25.333333333333
$ws.Range("I25").Value = 372
$ws.Range("J25").Value = 326
$ws.Range("K25").Value = 14.110429447852
$ws.Range("L25").Value = 24.832214765100
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 39
$ws.Range("E26").Value = -41.025641025641
$ws.Range("F26").Value = 119
$ws.Range("G26").Value = 122
$ws.Range("H26").Value = -2.459016393442
$ws.Range("I26").Value = 516
$ws.Range("J26").Value = 518
$ws.Range("K26").Value = -0.386100386100
$ws.Range("L26").Value = 24.337349397590
$ws.Range("M26").Value = -17.834394904458
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 11
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 33
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = -17.5
$ws.Range("L27").Value = 13.793103448275
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -46.153846153846
$ws.Range("I28").Value = 46
$ws.Range("J28").Value = 44
$ws.Range("K28").Value = 4.545454545454
$ws.Range("L28").Value = 24.324324324324
$ws.Range("J29").Value = 18
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -41.935483870967
$ws.Range("M29").Value = -40
$ws.Range("N29").Value = -86.861313868613
$ws.Range("J30").Value = 17
$ws.Range("K30").Value = -23.529411764705
$ws.Range("L30").Value = -50
$ws.Range("M30").Value = -45.833333333333
$ws.Range("N30").Value = -89.84375

# --- Shared string text edits (header volume number and date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21,2).Text = "19"

$c9 = $ws.Range("C9")
$c9.Characters(47,8).Text = "5/11/2025"
$c9.Characters(27,9).Text = "5/5/2025"
